# Insert a newly-stocked item ("DEXAFLOX EYE DROPS 5 ML") into the low-stock
# report, renumber the following rows, refresh the running total, and bump
# the generated-at timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 currently holds "DICLAC 150 ID 20 PROLONGED R TABS" - the new
# product is alphabetically just before it, so insert a fresh row above it.
# Inserting here shifts rows 13..37 down to 14..38 and carries their
# formatting/merges along automatically.
$ws.Range("A13:Q13").Insert()

# New row 13: copy the look of the row right above it (row 12, the same
# product-row layout used throughout the table), then fill in the data
# for the new product.
$ws.Range("A12:Q12").Copy()
$ws.Range("A13:Q13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 3).Value = "DEXAFLOX EYE DROPS 5 ML"
$ws.Cells.Item(13, 8).Value = "1:0"

# Columns L and P carry a numeric-looking display format but actually
# store plain text in this report, so flip to text format while writing
# the value and then restore the original display format.
$fmtL = $ws.Cells.Item(13, 12).NumberFormat
$ws.Cells.Item(13, 12).NumberFormat = "@"
$ws.Cells.Item(13, 12).Value = "1"
$ws.Cells.Item(13, 12).NumberFormat = $fmtL

$ws.Cells.Item(13, 14).Value = "53.00"

$fmtP = $ws.Cells.Item(13, 16).NumberFormat
$ws.Cells.Item(13, 16).NumberFormat = "@"
$ws.Cells.Item(13, 16).Value = "53.0000"
$ws.Cells.Item(13, 16).NumberFormat = $fmtP

$ws.Cells.Item(13, 17).Value = "1:0"

# Renumber the sequential index column (A) for every product row so it
# stays 1..30 with no gaps after the insertion.
for ($r = 14; $r -le 36; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# The running total (now on row 37) needs to include the new product's
# sell price.
$ws.Cells.Item(37, 16).Value = 1893.29

# Footer timestamp (now on row 38) reflects the new export time.
$ws.Cells.Item(38, 1).Value = "Tuesday, 30 September, 2025 2:37 PM"
